$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Clear out the stale URL / UserName / Password block (U2:W2) that used to hold
# the Oracle login hyperlink + credentials, and drop the hyperlink that was
# attached to U2 along with it.
$ws.Activate()
$ws.Range("U2:W2").Select()
$ws.Range("U2:W2").ClearContents()
$ws.Hyperlinks.Delete()
